$d = $word.ActiveDocument
$d.Content.Find.Execute("Proof: Scalar product", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Proof: Scalar product", 2)
